$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

$data = @(
    @(44370, 12467, 179, 5873, 18519, 5397, 346, 23, 323, 130),
    @(44371, 12591, 131, 5887, 18609, 5430, 327, 23, 304, 130),
    @(44372, 12645, 109, 5914, 18668, 5448, 331, 20, 311, 135),
    @(44373, 12736, 81, 5934, 18751, 5499, 300, 20, 280, 135),
    @(44374, 12738, 79, 5934, 18751, 5545, 253, 20, 233, 136)
)

$startRow = 357
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($col = 1; $col -le 10; $col++) {
        $ws.Cells.Item($row, $col).Value = $rowData[$col - 1]
    }
}

$ws.Range("A358").Select()
